$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading-percent results for the 380 kV case (Case_1_48)
$ws.Range("B2").Value = 19.81818561421435
$ws.Range("D2").Value = 3.40469112619871
$ws.Range("E2").Value = 29.66130923871608
$ws.Range("F2").Value = 22.59753965086227
$ws.Range("G2").Value = 27.38955477987595
$ws.Range("H2").Value = 12.8946362081942
$ws.Range("L2").Value = 12.43222193316058
$ws.Range("B3").Value = 19.29738277604039
$ws.Range("D3").Value = 3.454813411382561
$ws.Range("E3").Value = 28.47873392376783
$ws.Range("F3").Value = 22.43358088938869
$ws.Range("G3").Value = 26.90633163319794
$ws.Range("H3").Value = 12.94907333789256
$ws.Range("L3").Value = 11.95749511381063
$ws.Range("B4").Value = 18.9720676584672
$ws.Range("D4").Value = 3.487086607094113
$ws.Range("E4").Value = 27.72841178618647
$ws.Range("F4").Value = 22.34926239252703
$ws.Range("G4").Value = 26.63284226854571
$ws.Range("H4").Value = 12.98922935257285
$ws.Range("L4").Value = 11.65612534047784
$ws.Range("B5").Value = 18.83827881030411
$ws.Range("D5").Value = 3.500615369930876
$ws.Range("E5").Value = 27.41697770047072
$ws.Range("F5").Value = 22.31902890305325
$ws.Range("G5").Value = 26.52742481061743
$ws.Range("H5").Value = 13.00726514075566
$ws.Range("L5").Value = 11.53097836979621
$ws.Range("B6").Value = 18.81599481633678
$ws.Range("D6").Value = 3.502884607632663
$ws.Range("E6").Value = 27.3649349628025
$ws.Range("F6").Value = 22.31425809082296
$ws.Range("G6").Value = 26.51028975793055
$ws.Range("H6").Value = 13.01036032448364
$ws.Range("L6").Value = 11.51006141492743
$ws.Range("B7").Value = 18.97026804053109
$ws.Range("D7").Value = 3.48726753280337
$ws.Range("E7").Value = 27.72423405231213
$ws.Range("F7").Value = 22.34883793309862
$ws.Range("G7").Value = 26.63139591297674
$ws.Range("H7").Value = 12.98946584810733
$ws.Range("L7").Value = 11.65444681403473
$ws.Range("B8").Value = 19.63986314910195
$ws.Range("D8").Value = 3.421662685088204
$ws.Range("E8").Value = 29.2588182590984
$ws.Range("F8").Value = 22.53762567331855
$ws.Range("G8").Value = 27.21824728849541
$ws.Range("H8").Value = 12.91199678476945
$ws.Range("L8").Value = 12.27067218990298
$ws.Range("B9").Value = 20.90192000528388
$ws.Range("D9").Value = 3.30487696583143
$ws.Range("E9").Value = 32.06014332775068
$ws.Range("F9").Value = 23.03627513159852
$ws.Range("G9").Value = 28.54336564170821
$ws.Range("H9").Value = 12.8143762006081
$ws.Range("L9").Value = 13.39494665858717
$ws.Range("B10").Value = 21.789381463241
$ws.Range("D10").Value = 3.226285454446914
$ws.Range("E10").Value = 33.97346305386724
$ws.Range("F10").Value = 23.4782738122234
$ws.Range("G10").Value = 29.60858305725608
$ws.Range("H10").Value = 12.77689444673939
$ws.Range("L10").Value = 14.16313150691924
$ws.Range("B11").Value = 22.18286811910625
$ws.Range("D11").Value = 3.192095721532263
$ws.Range("E11").Value = 34.80944764243413
$ws.Range("F11").Value = 23.69496923578401
$ws.Range("G11").Value = 30.10973084701969
$ws.Range("H11").Value = 12.76749786614496
$ws.Range("L11").Value = 14.49895109941989
$ws.Range("B12").Value = 22.33028229028006
$ws.Range("D12").Value = 3.179373702217696
$ws.Range("E12").Value = 35.12086509375368
$ws.Range("F12").Value = 23.77919188827813
$ws.Range("G12").Value = 30.30160185692709
$ws.Range("H12").Value = 12.76505636153309
$ws.Range("L12").Value = 14.62408241897254
$ws.Range("B13").Value = 22.29860655798972
$ws.Range("D13").Value = 3.182103608280733
$ws.Range("E13").Value = 35.05402785325231
$ws.Range("F13").Value = 23.76095815264069
$ws.Range("G13").Value = 30.26019025937583
$ws.Range("H13").Value = 12.76553227281804
$ws.Range("L13").Value = 14.59722480914326
$ws.Range("B14").Value = 22.19502850220156
$ws.Range("D14").Value = 3.191044563982689
$ws.Range("E14").Value = 34.83517246291527
$ws.Range("F14").Value = 23.70185539976843
$ws.Range("G14").Value = 30.12547579056278
$ws.Range("H14").Value = 12.76727453297978
$ws.Range("L14").Value = 14.50928693502344
$ws.Range("B15").Value = 22.13137341130666
$ws.Range("D15").Value = 3.196550456196189
$ws.Range("E15").Value = 34.70044072021328
$ws.Range("F15").Value = 23.66593258944826
$ws.Range("G15").Value = 30.04322394885791
$ws.Range("H15").Value = 12.76848760902199
$ws.Range("L15").Value = 14.45515519570118
$ws.Range("B16").Value = 21.76344989826585
$ws.Range("D16").Value = 3.22855128477479
$ws.Range("E16").Value = 33.91811907758635
$ws.Range("F16").Value = 23.46441987237533
$ws.Range("G16").Value = 29.57614089824757
$ws.Range("H16").Value = 12.77766396626772
$ws.Range("L16").Value = 14.14090384732088
$ws.Range("B17").Value = 21.53503285783881
$ws.Range("D17").Value = 3.248582958428345
$ws.Range("E17").Value = 33.42922433260612
$ws.Range("F17").Value = 23.34474507866216
$ws.Range("G17").Value = 29.29364184755022
$ws.Range("H17").Value = 12.78526568028849
$ws.Range("L17").Value = 13.94457174610987
$ws.Range("B18").Value = 21.40269757041039
$ws.Range("D18").Value = 3.260251637685924
$ws.Range("E18").Value = 33.14479785249047
$ws.Range("F18").Value = 23.27738780634647
$ws.Range("G18").Value = 29.13273318630036
$ws.Range("H18").Value = 12.79035754375469
$ws.Range("L18").Value = 13.83036744229516
$ws.Range("B19").Value = 21.35773080192313
$ws.Range("D19").Value = 3.264227691647097
$ws.Range("E19").Value = 33.04794845804581
$ws.Range("F19").Value = 23.25483775648342
$ws.Range("G19").Value = 29.07853192762685
$ws.Range("H19").Value = 12.79220463112074
$ws.Range("L19").Value = 13.79148261308664
$ws.Range("B20").Value = 21.55944807632003
$ws.Range("D20").Value = 3.246435338087229
$ws.Range("E20").Value = 33.48160360155494
$ws.Range("F20").Value = 23.35733242330853
$ws.Range("G20").Value = 29.32355318623971
$ws.Range("H20").Value = 12.78438187187856
$ws.Range("L20").Value = 13.96560461954419
$ws.Range("B21").Value = 22.22549596404688
$ws.Range("D21").Value = 3.188412281987718
$ws.Range("E21").Value = 34.8995969024213
$ws.Range("F21").Value = 23.71915725673366
$ws.Range("G21").Value = 30.16499004759664
$ws.Range("H21").Value = 12.76673235854323
$ws.Range("L21").Value = 14.53517223596682
$ws.Range("B22").Value = 22.65147686429621
$ws.Range("D22").Value = 3.15180194885151
$ws.Range("E22").Value = 35.79623451195182
$ws.Range("F22").Value = 23.96820508838115
$ws.Range("G22").Value = 30.72701159340314
$ws.Range("H22").Value = 12.76171325055488
$ws.Range("L22").Value = 14.89552086890625
$ws.Range("B23").Value = 22.42501257589694
$ws.Range("D23").Value = 3.171221469585658
$ws.Range("E23").Value = 35.32049681001573
$ws.Range("F23").Value = 23.83416230680858
$ws.Range("G23").Value = 30.4260373068973
$ws.Range("H23").Value = 12.76379088928704
$ws.Range("L23").Value = 14.7043069203792
$ws.Range("B24").Value = 21.54841310918005
$ws.Range("D24").Value = 3.247405803183909
$ws.Range("E24").Value = 33.45793336625693
$ws.Range("F24").Value = 23.3516371766549
$ws.Range("G24").Value = 29.31002556424174
$ws.Range("H24").Value = 12.78477919547466
$ws.Range("L24").Value = 13.95609979536262
$ws.Range("B25").Value = 20.56690110660978
$ws.Range("D25").Value = 3.335203121139938
$ws.Range("E25").Value = 31.32660251772367
$ws.Range("F25").Value = 22.88785489944406
$ws.Range("G25").Value = 28.16767716871809
$ws.Range("H25").Value = 12.83485046191471
$ws.Range("L25").Value = 13.10053965242146
